# Apply the changes described by the commit to the "Casos de Uso" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")
$ws.Activate()

# Row 10 (Tarea "Asegurar el registro y la consulta de imágenes en Linux."):
#   - Estatus (F10) goes from "Por iniciar" to "Hecho"
#   - Día 4 "Cons." (Q10) gets 3 hours consumed registered
$ws.Range("F10").Value = "Hecho"
$ws.Range("Q10").Value = 3

# Row 12 (Tarea "Instalador de aplicación para Linux."):
#   - Día 4 "Cons." (Q12) gets 1 hour consumed registered
$ws.Range("Q12").Value = 1

# Update the frozen-pane anchor and the remembered selection in the
# bottom-right pane to match where the user was working afterwards.
$ws.Activate()
$ws.Range("G9").Select()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("G6").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("D15").Select()

$wb.Save()
